$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.613015413284302
$ws.Range("B1").Value = 6.223897933959961
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.597001552581787
$ws.Range("E1").Value = 1.991630434989929
